$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- 1. Fill in the previously-empty rows 9 and 10 with real data ---
$ws.Range("B9").Value = 0.44444444444444442
$ws.Range("C9").Value = 0.51041666666666663
$ws.Range("E9").Value = "User stories"

$ws.Range("B10").Value = 0.63888888888888895
$ws.Range("C10").Value = 0.67361111111111116
$ws.Range("E10").Value = "meeting avec chef de projet"
$ws.Range("F10").Value = "Création du sprint 1 et ajout des user stories + revue"
$ws.Range("F9").Value = "Création des user stories dans ice scrum"

# --- 2. Fix capitalisation / wording on a few existing "Description" cells ---
$ws.Range("F8").Value = "Création du mcd su draw.io"
$ws.Range("F6").Value = "PhpStorm, MySQL workbench, HeidiSQL"
$ws.Range("F4").Value = "Intro, objectif ajouter dans la doc"

# --- 3. Add new "Type" column (G) to the table, extending it from A1:F11 to A1:G11 ---
$null = $tbl.ListColumns.Add()

# Header text for the new column
$ws.Range("G1").Value = "Type"

# Copy formatting (number format / style) from existing cells so the new
# column's cells carry the same style indices as their row neighbours
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("G2:G11").PasteSpecial(-4122)

# --- 4. Update the selected cell (cosmetic, matches the author's last selection) ---
$ws.Range("E17").Select()
